# "Improve CT in Study Design sheet"
#
# The studyDesignBlindingScheme / trialIntentTypes / trialTypes /
# interventionModel controlled-terminology values on the "studyDesign"
# sheet are updated to drop the raw NCI-code prefixes (and, for
# trialTypes/interventionModel, to pick different coded terms), and the
# "studyDesign" sheet is left as the active/selected sheet with
# A1:E6 selected.

$wb = $excel.ActiveWorkbook

$studyDesign = $wb.Worksheets.Item("studyDesign")

# Update the controlled-terminology values (column B, merged B:E) for
# rows 3-6: studyDesignBlindingScheme, trialIntentTypes, trialTypes,
# interventionModel.
$studyDesign.Range("B3").Value = "OPEN LABEL"
$studyDesign.Range("B4").Value = "BASIC SCIENCE,    DEVICE FEASIBILITY"
$studyDesign.Range("B5").Value = "Efficacy Study"
$studyDesign.Range("B6").Value = "C82639"

# Make "studyDesign" the active sheet/tab, with A1:E6 selected, matching
# the author's final view state when the workbook was saved.
$studyDesign.Activate()
$studyDesign.Range("A1:E6").Select()
